$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("K2").Value = 4174
$ws.Range("K3").Value = 4272
$ws.Range("K4").Value = 861
$ws.Range("K5").Value = 309
$ws.Range("K6").Value = 4781
$ws.Range("K7").Value = 14397

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("K2").Value = 123
$ws.Range("K4").Value = 51
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 421
$ws.Range("K8").Value = 970
$ws.Range("K9").Value = 60
$ws.Range("K11").Value = 281
$ws.Range("K14").Value = 81
$ws.Range("K16").Value = 45
$ws.Range("K18").Value = 100
$ws.Range("K19").Value = 437
$ws.Range("K20").Value = 325
$ws.Range("K25").Value = 67
$ws.Range("K29").Value = 759
$ws.Range("K30").Value = 49
$ws.Range("K31").Value = 155
$ws.Range("K33").Value = 600
$ws.Range("K37").Value = 488
$ws.Range("K42").Value = 526
$ws.Range("K47").Value = 84
$ws.Range("K48").Value = 184
$ws.Range("K52").Value = 391
$ws.Range("K54").Value = 266
$ws.Range("K59").Value = 25
$ws.Range("K63").Value = 44
$ws.Range("K64").Value = 89
$ws.Range("K65").Value = 332
$ws.Range("K67").Value = 562
$ws.Range("K71").Value = 44
$ws.Range("K77").Value = 103
$ws.Range("K79").Value = 371
$ws.Range("K83").Value = 306
$ws.Range("K84").Value = 105
$ws.Range("K85").Value = 647
$ws.Range("K86").Value = 97
$ws.Range("K88").Value = 167
$ws.Range("K89").Value = 202
$ws.Range("K91").Value = 160
$ws.Range("K92").Value = 53
$ws.Range("K94").Value = 181
$ws.Range("K95").Value = 246
$ws.Range("K96").Value = 161
$ws.Range("K99").Value = 244
$ws.Range("K101").Value = 14397

# Sheet 3: Bridgeport
$ws = $wb.Worksheets.Item(3)
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 81

# Sheet 4: West Ridge
$ws = $wb.Worksheets.Item(4)
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 161

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item(5)
$ws.Range("K2").Value = 151
$ws.Range("K7").Value = 421

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 281

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item(7)
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 202

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item(8)
$ws.Range("K3").Value = 217
$ws.Range("K6").Value = 151
$ws.Range("K7").Value = 647

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item(9)
$ws.Range("K3").Value = 102
$ws.Range("K6").Value = 151
$ws.Range("K7").Value = 391

# Sheet 12: Austin
$ws = $wb.Worksheets.Item(12)
$ws.Range("K2").Value = 275
$ws.Range("K6").Value = 325
$ws.Range("K7").Value = 970

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item(13)
$ws.Range("K3").Value = 106
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 306

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item(14)
$ws.Range("K2").Value = 165
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 600

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item(15)
$ws.Range("K2").Value = 80
$ws.Range("K4").Value = 13
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 246

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item(16)
$ws.Range("K3").Value = 162
$ws.Range("K6").Value = 147
$ws.Range("K7").Value = 488

# Sheet 17: New City
$ws = $wb.Worksheets.Item(17)
$ws.Range("K3").Value = 85
$ws.Range("K6").Value = 132
$ws.Range("K7").Value = 332

# Sheet 18: Woodlawn
$ws = $wb.Worksheets.Item(18)
$ws.Range("K3").Value = 100
$ws.Range("K7").Value = 244

# Sheet 19: Fuller Park
$ws = $wb.Worksheets.Item(19)
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 49

# Sheet 20: Gage Park
$ws = $wb.Worksheets.Item(20)
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 155

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item(21)
$ws.Range("K2").Value = 164
$ws.Range("K3").Value = 196
$ws.Range("K4").Value = 29
$ws.Range("K7").Value = 562

# Sheet 22: South Deering
$ws = $wb.Worksheets.Item(22)
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 105

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Range("K6").Value = 132
$ws.Range("K7").Value = 266

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Range("K2").Value = 217
$ws.Range("K3").Value = 270
$ws.Range("K6").Value = 210
$ws.Range("K7").Value = 759

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item(26)
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 184

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item(27)
$ws.Range("K6").Value = 133
$ws.Range("K7").Value = 437

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item(30)
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 106

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Range("K2").Value = 144
$ws.Range("K3").Value = 167
$ws.Range("K6").Value = 191
$ws.Range("K7").Value = 526

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item(40)
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 160

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Range("K3").Value = 121
$ws.Range("K5").Value = 15
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 371

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item(43)
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 89

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item(44)
$ws.Range("K2").Value = 111
$ws.Range("K3").Value = 100
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 325

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 100

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 181

# Sheet 52: East Side
$ws = $wb.Worksheets.Item(52)
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 67

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item(53)
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 84

# Sheet 61: Avalon Park
$ws = $wb.Worksheets.Item(61)
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 60

# Sheet 63: Montclare
$ws = $wb.Worksheets.Item(63)
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 25

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item(64)
$ws.Range("K2").Value = 37
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 123

# Sheet 66: West Elsdon
$ws = $wb.Worksheets.Item(66)
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 53

# Sheet 68: United Center
$ws = $wb.Worksheets.Item(68)
$ws.Range("K3").Value = 52
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 167

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item(72)
$ws.Range("K4").Value = 37
$ws.Range("K7").Value = 97

# Sheet 81: Oakland
$ws = $wb.Worksheets.Item(81)
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 44

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item(84)
$ws.Range("K2").Value = 46
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 103

# Sheet 90: Archer Heights
$ws = $wb.Worksheets.Item(90)
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 51

# Sheet 94: Bucktown
$ws = $wb.Worksheets.Item(94)
$ws.Range("K2").Value = 13
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 45
